$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.997.55"
$ws.Range("E2").Value = "  -2.38%  "
$ws.Range("D3").Value = "3.447.50"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'584.24"
$ws.Range("E5").Value = "  -1.19%  "
$ws.Range("D6").Value = "'174.02"
$ws.Range("E6").Value = "  -3.03%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.602"
$ws.Range("E8").Value = "  -1.79%  "
$ws.Range("D9").Value = "3.445.29"
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("E10").Value = "  -4.74%  "
$ws.Range("D11").Value = "'6.90"
$ws.Range("E11").Value = "  -1.30%  "
$ws.Range("E12").Value = "  -3.75%  "
$ws.Range("D13").Value = "4.045.22"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("D15").Value = "'29.08"
$ws.Range("E15").Value = "  -9.73%  "
$ws.Range("D16").Value = "66.038.54"
$ws.Range("E16").Value = "  -2.31%  "
$ws.Range("E17").Value = "  -3.02%  "
$ws.Range("D18").Value = "3.446.73"
$ws.Range("E18").Value = "  -0.63%  "
$ws.Range("E19").Value = "  -3.16%  "
$ws.Range("E20").Value = "  -1.34%  "
$ws.Range("D21").Value = "'367.78"
$ws.Range("E21").Value = "  -4.73%  "
$ws.Range("E22").Value = "  -2.78%  "
$ws.Range("D23").Value = "'72.76"
$ws.Range("E23").Value = "  +0.78%  "
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").Value = "'9.75"
$ws.Range("E27").Value = "  -3.23%  "
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  -2.58%  "
$ws.Range("E31").Value = "  -2.91%  "
$ws.Range("D32").Value = "'5.72"
$ws.Range("E32").Value = "  -5.05%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  -6.00%  "
$ws.Range("E35").Value = "  -3.58%  "
$ws.Range("E36").Value = "  -1.70%  "
$ws.Range("D37").Value = "'160.70"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").Value = "'28.87"
$ws.Range("E38").Value = "  +3.91%  "
$ws.Range("E39").Value = "  -0.89%  "
$ws.Range("D40").Value = "'2.64"
$ws.Range("E40").Value = "  -1.41%  "
$ws.Range("E41").Value = "  -4.77%  "
$ws.Range("D42").Value = "2.764.23"
$ws.Range("E42").Value = "  +1.52%  "
$ws.Range("D43").Value = "'4.46"
$ws.Range("E43").Value = "  -1.73%  "
$ws.Range("E44").Value = "  -2.86%  "
$ws.Range("D45").Value = "'0.0681"
$ws.Range("E45").Value = "  -4.01%  "
$ws.Range("D46").Value = "'40.14"
$ws.Range("E46").Value = "  -3.35%  "
$ws.Range("D47").Value = "'24.35"
$ws.Range("E47").Value = "  -5.76%  "
$ws.Range("D48").Value = "'0.0290"
$ws.Range("E48").Value = "  -2.43%  "
$ws.Range("D49").Value = "'325.13"
$ws.Range("E49").Value = "  -0.90%  "
$ws.Range("D50").Value = "'6.24"
$ws.Range("E50").Value = "  -0.65%  "
$ws.Range("E51").Value = "  -3.31%  "
